$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45190 -> 2023-09-21)
# that was bumped to 45192 (-> 2023-09-23) for every data row (rows 2-260).
$ws.Range("C2:C260").Value = 45192
